$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.556.65'
$ws.Range("E2").Value = '  +2.08%  '
$ws.Range("D3").Value = '3.119.08'
$ws.Range("E3").Value = '  +2.27%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '394.95'
$ws.Range("E5").Value = '  +2.73%  '
$ws.Range("D6").Value = '103.66'
$ws.Range("E6").Value = '  +1.01%  '
$ws.Range("E7").Value = '  -0.47%  '
$ws.Range("D9").Value = '0.604'
$ws.Range("E9").Value = '  +3.45%  '
$ws.Range("D10").Value = '38.08'
$ws.Range("E10").Value = '  +3.21%  '
$ws.Range("E11").Value = '  +0.96%  '
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("D13").Value = '3.601.76'
$ws.Range("E13").Value = '  +1.96%  '
$ws.Range("D14").Value = '18.79'
$ws.Range("E14").Value = '  +0.67%  '
$ws.Range("D15").Value = '7.83'
$ws.Range("E15").Value = '  +1.27%  '
$ws.Range("D16").Value = '1.05'
$ws.Range("E16").Value = '  +8.00%  '
$ws.Range("D17").Value = '3.117.54'
$ws.Range("E17").Value = '  +1.90%  '
$ws.Range("D18").Value = '10.74'
$ws.Range("E18").Value = '  +2.60%  '
$ws.Range("D19").Value = '52.471.64'
$ws.Range("E19").Value = '  +1.76%  '
$ws.Range("E20").Value = '  +3.47%  '
$ws.Range("D21").Value = '12.76'
$ws.Range("E21").Value = '  +2.98%  '
$ws.Range("E22").Value = '  +0.86%  '
$ws.Range("D23").Value = '70.98'
$ws.Range("E23").Value = '  +1.12%  '
$ws.Range("D24").Value = '268.83'
$ws.Range("E24").Value = '  +0.30%  '
$ws.Range("E25").Value = '  +1.82%  '
$ws.Range("D26").Value = '8.03'
$ws.Range("E26").Value = '  -3.80%  '
$ws.Range("D27").Value = '27.54'
$ws.Range("E27").Value = '  +2.06%  '
$ws.Range("D28").Value = '7.45'
$ws.Range("E28").Value = '  +2.98%  '
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("E30").Value = '  -2.49%  '
$ws.Range("E31").Value = '  +0.10%  '
$ws.Range("D32").Value = '10.89'
$ws.Range("E32").Value = '  +6.15%  '
$ws.Range("B33").Value = 'InjectiveProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D33").Value = '36.72'
$ws.Range("E33").Value = '  +6.41%  '
$ws.Range("B34").Value = 'VeChain'
$ws.Range("C34").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D34").Value = '0.0489'
$ws.Range("E34").Value = '  +10.53%  '
$ws.Range("E35").Value = '  +0.96%  '
$ws.Range("D36").Value = '50.06'
$ws.Range("E36").Value = '  -0.83%  '
$ws.Range("E38").Value = '  +1.44%  '
$ws.Range("E39").Value = '  +10.56%  '
$ws.Range("E40").Value = '  +1.13%  '
$ws.Range("E41").Value = '  +6.26%  '
$ws.Range("D42").Value = '130.27'
$ws.Range("E42").Value = '  +1.65%  '
$ws.Range("D43").Value = '16.98'
$ws.Range("E43").Value = '  +0.18%  '
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("E45").Value = '  +0.38%  '
$ws.Range("D46").Value = '22.11'
$ws.Range("E46").Value = '  +0.74%  '
$ws.Range("D47").Value = '2.45'
$ws.Range("E47").Value = '  -2.67%  '
$ws.Range("E48").Value = '  -0.72%  '
$ws.Range("D49").Value = '2.079.77'
$ws.Range("E49").Value = '  +2.03%  '
$ws.Range("D50").Value = '0.0537'
$ws.Range("E50").Value = '  +37.45%  '
$ws.Range("D51").Value = '0.916'
$ws.Range("E51").Value = '  +10.11%  '
